# Supermere plasma vs Control+LPS heatmap
#
# 1) The "Isolation" column (F) used the shared text "NA" for every
#    untreated/control sample. Re-label every occurrence as "None".
# 2) Move the sheet's selection from the old last-row cell (F90) to the
#    full data range of column F (F2:F89), reflecting the new heatmap
#    focus on the Isolation column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$lastCol = $used.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value() -eq "NA") {
            $cell.Value = "None"
        }
    }
}

# Reselect so the saved view highlights the Isolation column's data range.
$ws.Range("F2:F89").Select()
